$d = $word.ActiveDocument

$find = $d.Content.Find
$find.ClearFormatting()
$find.Replacement.ClearFormatting()
$find.Execute(
    "Datas das campanhas de Constelação de Gêmeos 2022:",
    $true,
    $false,
    $false,
    $false,
    $false,
    $true,
    1,
    $false,
    "Datas das campanhas de 2022 que usam Constelação de Gêmeos:",
    2
)
